$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 19:22"

# Update country statistics cells (Casos totales / Nuevos casos / Casos activos /
# Recuperados / Casos criticos / Muertes hoy / Muertes) per the source diff.
$ws.Range("B4").Value = 385838
$ws.Range("C4").Value = 18834
$ws.Range("D4").Value = 21311
$ws.Range("E4").Value = 352300
$ws.Range("G4").Value = 1356
$ws.Range("H4").Value = 12227
$ws.Range("B13").Value = 22253
$ws.Range("C13").Value = 596
$ws.Range("D13").Value = 8704
$ws.Range("E13").Value = 12728
$ws.Range("B17").Value = 12613
$ws.Range("C17").Value = 316
$ws.Range("E17").Value = 8324
$ws.Range("E25").Value = 5748
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 89
$ws.Range("B26").Value = 5709
$ws.Range("C26").Value = 345
$ws.Range("E26").Value = 5474
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = 210
$ws.Range("B27").Value = 5311
$ws.Range("C27").Value = 533
$ws.Range("E27").Value = 4740
$ws.Range("B31").Value = 4848
$ws.Range("C31").Value = 435
$ws.Range("E31").Value = 4528
$ws.Range("E52").Value = 1234
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 56
$ws.Range("B60").Value = 1184
$ws.Range("C60").Value = 64
$ws.Range("D60").Value = 93
$ws.Range("E60").Value = 1001
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 10
$ws.Range("H60").Value = 90
$ws.Range("B61").Value = 1160
$ws.Range("C61").Value = 54
$ws.Range("D61").Value = 241
$ws.Range("E61").Value = 918
$ws.Range("F61").Value = 14
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 1
$ws.Range("B62").Value = 1149
$ws.Range("C62").Value = 41
$ws.Range("D62").Value = 69
$ws.Range("E62").Value = 1059
$ws.Range("F62").Value = 12
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 21
$ws.Range("B72").Value = 764
$ws.Range("C72").Value = 90
$ws.Range("E72").Value = 663
$ws.Range("B95").Value = 353
$ws.Range("C95").Value = 4
$ws.Range("D95").Value = 138
$ws.Range("E95").Value = 209
$ws.Range("F95").Value = 5
$ws.Range("H95").Value = 6
$ws.Range("D96").Value = 40
$ws.Range("E96").Value = 309
$ws.Range("F96").Value = 4
$ws.Range("H96").Value = 0
$ws.Range("B131").Value = 79
$ws.Range("C131").Value = 2
$ws.Range("D131").Value = 4
$ws.Range("E131").Value = 74
$ws.Range("H131").Value = 1
$ws.Range("B132").Value = 78
$ws.Range("C132").Value = 9
$ws.Range("D132").Value = 5
$ws.Range("E132").Value = 69
$ws.Range("H132").Value = 4
$ws.Range("D178").Value = 7
$ws.Range("E178").Value = 5
$ws.Range("C192").Value = 1
$ws.Range("D192").Value = 1
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0
$ws.Range("B193").Value = 8
$ws.Range("C193").Value = 3
$ws.Range("F193").Value = 1
$ws.Range("G193").Value = 1
$ws.Range("H193").Value = 1
$ws.Range("E194").Value = 7
$ws.Range("F194").Value = 0
$ws.Range("H194").Value = 0
$ws.Range("D195").Value = 0
$ws.Range("F195").Value = 1
$ws.Range("H195").Value = 1
